# "update plots for each sample" - re-run of the allele-calling pipeline.
# Four wildtype peaks that previously failed detection (min_height too high
# relative to the real peak height) now clear the (lowered) min_height
# threshold, so they get detected: peak/size/height/status get populated,
# the failure message clears, the marker genotype/phenotype gets called,
# and the overall sample genotype is resolved.

$wb = $excel.ActiveWorkbook

$peak   = $wb.Worksheets.Item("peak_table")
$allele = $wb.Worksheets.Item("allele_table")
$marker = $wb.Worksheets.Item("marker_table")
$geno   = $wb.Worksheets.Item("genotype_result")

# --- peak_table: lower w_height (min-height threshold) for the four wildtype markers ---
# row 2  -> S1 / CYP2D6_14
# row 4  -> S1 / CYP2D6_49
# row 12 -> S2 / CYP2D6_4
# row 14 -> S2 / CYP2D6_17
$peak.Cells.Item(2, 14).Value  = 800
$peak.Cells.Item(4, 14).Value  = 400
$peak.Cells.Item(12, 14).Value = 400
$peak.Cells.Item(14, 14).Value = 700

# --- allele_table: same four wildtype rows now resolve to a detected peak ---
# columns: A sample, B gene, C marker, D label, E panel, F direction, G base,
#          H basetype, I min_bin, J max_bin, K min_height, L is_forward,
#          M is_detected, N peak, O size, P height, Q status, R message, S color

# row 2: CYP2D6_14 wildtype (S1, Forward)
$allele.Cells.Item(2, 11).Value = 800
$allele.Cells.Item(2, 13).Value = $true
$allele.Cells.Item(2, 14).Value = 36
$allele.Cells.Item(2, 15).Value = 29.11
$allele.Cells.Item(2, 16).Value = 949
$allele.Cells.Item(2, 17).Value = "ok"
$allele.Cells.Item(2, 18).Value = ""

# row 6: CYP2D6_49 wildtype (S1, Reverse)
$allele.Cells.Item(6, 11).Value = 400
$allele.Cells.Item(6, 13).Value = $true
$allele.Cells.Item(6, 14).Value = 17
$allele.Cells.Item(6, 15).Value = 38.82
$allele.Cells.Item(6, 16).Value = 471
$allele.Cells.Item(6, 17).Value = "ok"
$allele.Cells.Item(6, 18).Value = ""

# row 22: CYP2D6_4 wildtype (S2, Forward)
$allele.Cells.Item(22, 11).Value = 400
$allele.Cells.Item(22, 13).Value = $true
$allele.Cells.Item(22, 14).Value = 40
$allele.Cells.Item(22, 15).Value = 30.5
$allele.Cells.Item(22, 16).Value = 694
$allele.Cells.Item(22, 17).Value = "ok"
$allele.Cells.Item(22, 18).Value = ""

# row 26: CYP2D6_17 wildtype (S2, Forward)
$allele.Cells.Item(26, 11).Value = 700
$allele.Cells.Item(26, 13).Value = $true
$allele.Cells.Item(26, 14).Value = 26
$allele.Cells.Item(26, 15).Value = 38.77
$allele.Cells.Item(26, 16).Value = 787
$allele.Cells.Item(26, 17).Value = "ok"
$allele.Cells.Item(26, 18).Value = ""

# --- marker_table: genotype/phenotype now callable for the four markers ---
# columns: A sample, B gene, C marker, D label, E panel, F direction, G genotype, H phenotype
$marker.Cells.Item(2, 7).Value  = "GG"
$marker.Cells.Item(2, 8).Value  = "wildtype"

$marker.Cells.Item(4, 7).Value  = "TT"
$marker.Cells.Item(4, 8).Value  = "wildtype"

$marker.Cells.Item(12, 7).Value = "GG"
$marker.Cells.Item(12, 8).Value = "wildtype"

$marker.Cells.Item(14, 7).Value = "CC"
$marker.Cells.Item(14, 8).Value = "wildtype"

# --- genotype_result: overall sample genotype call ---
$geno.Cells.Item(2, 2).Value = "*1/*1"
